$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Add a new "valid credentials" scenario row (row 3), mirroring the
# existing row 2 data, and fill in row 4 (previously an empty, styled
# placeholder cell) with the same valid credentials.
$ws.Range("A3").Value = "standard_user"
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("A3").Style = $ws.Range("A2").Style

$ws.Range("A4").Value = "standard_user"
$ws.Range("B4").Value = "secret_sauce"
$ws.Range("A4").Style = $ws.Range("A2").Style

# Leave the selection where the editing session ended.
$ws.Range("E7").Select() | Out-Null
